# "reverts to absolute time specification for windowing, removes unneeded
# columns in .xlsx files"
#
# The "Start (min since recording start time)" column (column B) is no
# longer needed, so delete it outright. The remaining columns
# ("Elapsed Time (min)" and "# Blocks to Divide into") shift left into
# columns B and C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the (now-trailing, empty) column D the way the original editor's
# window selection ended up, then remove column B, shifting C:D to B:C.
$ws.Range("D1:D13").Select() | Out-Null
$ws.Columns.Item(2).Delete()
